$wb = $excel.ActiveWorkbook

# --- 1. Add the new "2023-2" worksheet after the last existing sheet ("2024-1") ---
$srcSheet  = $wb.Worksheets.Item("2024-1")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet  = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "2023-2"

# Copy the header row (A1:E1) plus the first two data rows (A2:E3) from
# "2024-1" so the new sheet starts out with identical layout/formatting.
$srcSheet.Range("A1:E3").Copy($newSheet.Range("A1:E3"))
# Columns F/G only have header content on row 1 for this sheet.
$srcSheet.Range("F1").Copy($newSheet.Range("F1"))
$srcSheet.Range("G1").Copy($newSheet.Range("G1"))
# F2 carries the row-2 style (no value) - Copy() drops style-only cells, so
# restore its formatting explicitly (cellXfs 9 == the "Normal 3" cell style).
$newSheet.Cells.Item(2, 6).Style = "Normal 3"

# --- 2. Fill in the "2023-2" season's data (same vessel/contract as the
#        "A.S/0001" rows elsewhere in the workbook, re-keyed to this season) ---
# NOTE: setting .Value resets a cell's direct style, so re-apply "Normal 3"
# (== style index 9) afterwards on every cell that needs it.
$newSheet.Cells.Item(2, 1).Value = "2023-2"
$newSheet.Cells.Item(2, 2).Value = "E/P VELA I"
$newSheet.Cells.Item(2, 2).Style = "Normal 3"
$newSheet.Cells.Item(2, 4).Value = "A.S/0001"
$newSheet.Cells.Item(2, 4).Style = "Normal 3"
$newSheet.Cells.Item(2, 5).Value = "A.S/0001-223"
$newSheet.Cells.Item(2, 5).Style = "Normal 3"

$newSheet.Cells.Item(3, 1).Value = "2023-2"
$newSheet.Cells.Item(3, 2).Value = "E/P VELA I"
$newSheet.Cells.Item(3, 2).Style = "Normal 3"
$newSheet.Cells.Item(3, 4).Value = "A.S/0001"
$newSheet.Cells.Item(3, 4).Style = "Normal 3"
$newSheet.Cells.Item(3, 5).Value = "A.S/0001-232"
$newSheet.Cells.Item(3, 5).Style = "Normal 3"

# --- 3. Column widths for the new sheet (C and E are the only custom ones) ---
$newSheet.Columns.Item(3).ColumnWidth = 23.26
$newSheet.Columns.Item(5).ColumnWidth = 11.59

# --- 4. Selection on the "2024-1" sheet moves to A2:G2 ---
[void]$srcSheet.Range("A2:G2").Select()

# --- 5. New sheet becomes the active tab, with its own selection at C7 ---
$newSheet.Activate()
[void]$newSheet.Range("C7").Select()
